$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-9, columns B-G (B:F are error metrics, G is a count)
$data = @{
    2 = @(0.0646005966593584, 0.4743505096951567, 0.5536090028381865, 0.7440490594296767, 0.7595443068293734, 21)
    3 = @(0.8918866930342709, 0.9430961204109612, 2.371025881809898, 1.539813586707786, 1.287823136188919, 20)
    4 = @(0.4882202679668279, 1.18756317707376, 3.637589225426803, 1.907246503582272, 1.894221823688656, 19)
    5 = @(0.5864295595391767, 0.7880825022580396, 0.8304769580234592, 0.9113050850420288, 0.7177740073227471, 18)
    6 = @(0.659487078565832, 0.7999889958063582, 1.154946370456449, 1.074684311998853, 0.8746568535088345, 17)
    7 = @(0.4265915177156135, 0.5072447508993148, 0.3397562009063173, 0.582886095996737, 0.4102368459569428, 16)
    8 = @(0.2793977770427775, 0.4649990375520733, 0.3217702627042455, 0.5672479728515964, 0.5109939315628624, 15)
    9 = @(0.4357956396000112, 0.4357956396000112, 0.24508204164741, 0.4950576144727096, 0.2463343710656916, 11)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}
